$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.817.61"
$ws.Range("E2").Value = "  -0.69%  "
$ws.Range("D3").Value = "2.213.66"
$ws.Range("E3").Value = "  -1.23%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "255.52"
$ws.Range("E5").Value = "  +3.03%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.618"
$ws.Range("E6").Value = "  +0.37%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "77.87"
$ws.Range("E7").Value = "  +3.59%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.595"
$ws.Range("E9").Value = "  -1.51%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "42.31"
$ws.Range("E10").Value = "  +3.33%  "
$ws.Range("E11").Value = "  -2.66%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "6.97"
$ws.Range("E12").Value = "  +0.94%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.102"
$ws.Range("E13").Value = "  +0.94%  "
$ws.Range("D14").Value = "2.533.14"
$ws.Range("E14").Value = "  -1.80%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "14.45"
$ws.Range("E15").Value = "  -1.15%  "
$ws.Range("D16").Value = "2.218.99"
$ws.Range("E16").Value = "  -1.16%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.780"
$ws.Range("E17").Value = "  -1.36%  "
$ws.Range("D18").Value = "42.775.89"
$ws.Range("E18").Value = "  -0.48%  "
$ws.Range("E19").Value = "  -1.89%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "71.03"
$ws.Range("E20").Value = "  -0.26%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.98"
$ws.Range("E21").Value = "  -0.37%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.28"
$ws.Range("E22").Value = "  +3.38%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "229.96"
$ws.Range("E23").Value = "  +0.13%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.31"
$ws.Range("E24").Value = "  -4.94%  "
$ws.Range("E25").Value = "  -0.26%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "42.35"
$ws.Range("E26").Value = "  +7.69%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.78"
$ws.Range("E27").Value = "  -0.31%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "3.35"
$ws.Range("E28").Value = "  -2.55%  "
$ws.Range("E29").Value = "  -0.09%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.20"
$ws.Range("E30").Value = "  -2.17%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "173.91"
$ws.Range("E31").Value = "  +1.30%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "20.44"
$ws.Range("E32").Value = "  +0.97%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0868"
$ws.Range("E33").Value = "  +8.31%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.24"
$ws.Range("E34").Value = "  -0.95%  "
$ws.Range("E35").Value = "  -0.85%  "
$ws.Range("E36").Value = "  +8.02%  "
$ws.Range("E37").Value = "  -2.58%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.32"
$ws.Range("E38").Value = "  -3.06%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "13.12"
$ws.Range("E39").Value = "  +0.43%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.84"
$ws.Range("E40").Value = "  +17.39%  "
$ws.Range("E41").Value = "  -0.45%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.202"
$ws.Range("E42").Value = "  -1.74%  "
$ws.Range("B43").Value = "MultiversX"
$ws.Range("C43").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "60.83"
$ws.Range("E43").Value = "  +2.54%  "
$ws.Range("B44").Value = "THORChain"
$ws.Range("C44").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.32"
$ws.Range("E44").Value = "  -1.62%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.488"
$ws.Range("E45").Value = "  +1.12%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "103.46"
$ws.Range("E46").Value = "  -1.09%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "8.38"
$ws.Range("E47").Value = "  -3.55%  "
$ws.Range("E48").Value = "  -1.44%  "
$ws.Range("E49").Value = "  +0.54%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.13"
$ws.Range("E50").Value = "  -1.67%  "
$ws.Range("E51").Value = "  +19.23%  "
